$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.365.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.41%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.846.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9985"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'240.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.53%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.6307"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.36%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.07541"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.10%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.2957"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'24.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.41%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07702"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.11%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.850.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -3.60%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6832"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.74%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.00001001"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.15%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.114.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.00%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'29.401.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.56%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'227.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.53%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.14%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.541"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'157.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.66%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.1398"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.39%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'8.355"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'17.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.463"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.95%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.05674"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.98%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D33").Value = "'4.022"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.38%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.841"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.16%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.21%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7181"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.19%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.596"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.259.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.25%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.60%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.55%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.225"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.09%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.9083"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.31%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D44").Value = "'101.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.69%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'66.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'7.061"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.4037"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.100"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.61%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.83%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.682"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.1124"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.45%  "
$ws.Range("E51").Style = "Normal"
